# Sprint8.xlsx - "Updated progress of sprint, sprint8"
#
# The Sprint sheet's task-status dropdown cells (F5:F7) get filled in, and
# the "Implement fix" task (row 5) records 3 units of effort already spent
# on Day 5 (column K). Shared-string insertion order follows the order the
# cells are written in (To do, then Done, then In progress), matching the
# authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint")

$ws.Range("F7").Value = "To do"
$ws.Range("F5").Value = "Done"
$ws.Range("F6").Value = "In progress"

$ws.Range("K5").Value = 3

# Leave the selection where the author's last edit did.
$ws.Activate()
$ws.Range("F6").Select()
